$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.186238
$ws.Range("H2").Value = 0.558714
$ws.Range("I2").Value = 0.05023668284714279
$ws.Range("J2").Value = 0.05023668284714279
$ws.Range("M2").Value = 14.321881
$ws.Range("N2").Value = 42.965643
$ws.Range("O2").Value = 0.2949569176783066
$ws.Range("P2").Value = 0.2949569176783066
$ws.Range("Q2").Value = 2.667278473678
$ws.Range("R2").Value = 24.005506263102
$ws.Range("S2").Value = 0.01481765712697589
$ws.Range("T2").Value = 0.01481765712697589
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.186238
$ws.Range("H3").Value = 0.558714
$ws.Range("I3").Value = 0.05023668284714279
$ws.Range("J3").Value = 0.05023668284714279
$ws.Range("N3").Value = 81.25250700000001
$ws.Range("O3").Value = 0.557794259435499
$ws.Range("P3").Value = 0.557794259435499
$ws.Range("Q3").Value = 5.044101466222001
$ws.Range("R3").Value = 45.39691319599801
$ws.Range("S3").Value = 0.02802173330521805
$ws.Range("T3").Value = 0.02802173330521805
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.186238
$ws.Range("H4").Value = 0.558714
$ws.Range("I4").Value = 0.05023668284714279
$ws.Range("J4").Value = 0.05023668284714279
$ws.Range("M4").Value = 7.149790333333333
$ws.Range("N4").Value = 21.449371
$ws.Range("O4").Value = 0.1472488228861944
$ws.Range("P4").Value = 0.1472488228861943
$ws.Range("Q4").Value = 1.331562652099333
$ws.Range("R4").Value = 11.984063868894
$ws.Range("S4").Value = 0.007397292414948848
$ws.Range("T4").Value = 0.007397292414948846
$ws.Range("I5").Value = 0.659992587420158
$ws.Range("J5").Value = 0.6599925874201579
$ws.Range("M5").Value = 14.321881
$ws.Range("N5").Value = 42.965643
$ws.Range("O5").Value = 0.2949569176783066
$ws.Range("P5").Value = 0.2949569176783066
$ws.Range("Q5").Value = 35.041804542892
$ws.Range("R5").Value = 315.376240886028
$ws.Range("S5").Value = 0.1946693792759801
$ws.Range("T5").Value = 0.19466937927598
$ws.Range("I6").Value = 0.659992587420158
$ws.Range("J6").Value = 0.6599925874201579
$ws.Range("N6").Value = 81.25250700000001
$ws.Range("O6").Value = 0.557794259435499
$ws.Range("P6").Value = 0.557794259435499
$ws.Range("S6").Value = 0.3681400765329458
$ws.Range("T6").Value = 0.3681400765329458
$ws.Range("I7").Value = 0.659992587420158
$ws.Range("J7").Value = 0.6599925874201579
$ws.Range("M7").Value = 7.149790333333333
$ws.Range("N7").Value = 21.449371
$ws.Range("O7").Value = 0.1472488228861944
$ws.Range("P7").Value = 0.1472488228861943
$ws.Range("Q7").Value = 17.49362080185734
$ws.Range("R7").Value = 157.442587216716
$ws.Range("S7").Value = 0.097183131611232
$ws.Range("T7").Value = 0.09718313161123196
$ws.Range("G8").Value = 1.074241333333333
$ws.Range("H8").Value = 3.222724
$ws.Range("I8").Value = 0.2897707297326994
$ws.Range("J8").Value = 0.2897707297326994
$ws.Range("M8").Value = 14.321881
$ws.Range("N8").Value = 42.965643
$ws.Range("O8").Value = 0.2949569176783066
$ws.Range("P8").Value = 0.2949569176783066
$ws.Range("Q8").Value = 15.38515654128133
$ws.Range("R8").Value = 138.466408871532
$ws.Range("S8").Value = 0.08546988127535063
$ws.Range("T8").Value = 0.08546988127535063
$ws.Range("G9").Value = 1.074241333333333
$ws.Range("H9").Value = 3.222724
$ws.Range("I9").Value = 0.2897707297326994
$ws.Range("J9").Value = 0.2897707297326994
$ws.Range("N9").Value = 81.25250700000001
$ws.Range("O9").Value = 0.557794259435499
$ws.Range("P9").Value = 0.557794259435499
$ws.Range("Q9").Value = 29.09493381878534
$ws.Range("R9").Value = 261.854404369068
$ws.Range("S9").Value = 0.1616324495973352
$ws.Range("T9").Value = 0.1616324495973352
$ws.Range("G10").Value = 1.074241333333333
$ws.Range("H10").Value = 3.222724
$ws.Range("I10").Value = 0.2897707297326994
$ws.Range("J10").Value = 0.2897707297326994
$ws.Range("M10").Value = 7.149790333333333
$ws.Range("N10").Value = 21.449371
$ws.Range("O10").Value = 0.1472488228861944
$ws.Range("P10").Value = 0.1472488228861943
$ws.Range("Q10").Value = 7.680600300733778
$ws.Range("R10").Value = 69.12540270660401
$ws.Range("S10").Value = 0.04266839886001354
$ws.Range("T10").Value = 0.04266839886001354